$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1219
$ws.Range("F5").Value = 74
$ws.Range("F10").Value = 618
$ws.Range("F16").Value = 1006
$ws.Range("C18").Value = "上海·第一届妖妖动漫游戏展"
$ws.Range("E18").Value = "2024.05.02 10:00-05.04 17:00"
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 68
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=84642"
$ws.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202404/T1ytR8k81713936182881.jpeg"
$ws.Range("C19").Value = "上海·第五十八届燃梦星辰国潮嘉年华-随机宅舞"
$ws.Range("D19").Value = "周家嘴路3608号 宝龙旭辉广场"
$ws.Range("E19").Value = "2024.05.02 10:20-05.03 16:30"
$ws.Range("F19").Value = 707
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=82761"
$ws.Range("I19").Value = "//i0.hdslb.com/bfs/openplatform/202403/azEA4EM01710236719279.jpeg"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "2024-05-03"
$ws.Range("B20").ClearFormats()
$ws.Range("C20").Value = "上海·DizzyMart2024电则市场 中国同人音乐展会"
$ws.Range("D20").Value = "顾村镇蕰川路6号 智慧湾科创园"
$ws.Range("E20").Value = "2024.05.03 09:00-05.04 20:00"
$ws.Range("F20").Value = 657
$ws.Range("G20").Value = 138
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=84202"
$ws.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202404/etRgMvxv1712656961255.jpeg"
$ws.Range("C21").Value = "上海·HD动漫主题嘉年华·大唐夜市之剑侠奇缘"
$ws.Range("D21").Value = "天等路400号，与华东理工大学仅一墙之隔 品域凌云里"
$ws.Range("E21").Value = "2024.05.03 10:00-05.04 17:30"
$ws.Range("F21").Value = 57
$ws.Range("G21").Value = 75
$ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=84247"
$ws.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202403/xOubEDCI1711594090227.jpeg"
$ws.Range("C22").Value = "上海·坏孩纸物语第41届动漫节"
$ws.Range("D22").Value = "曹杨路1888号 复悦荟"
$ws.Range("E22").Value = "2024.05.03 12:00-05.04 16:00"
$ws.Range("F22").Value = 40
$ws.Range("G22").Value = 77.7
$ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=84369"
$ws.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202404/GeehIr1u1713248070279.png"
$ws.Range("F30").Value = 799
$ws.Range("F31").Value = 144
$ws.Range("F33").Value = 1329
$ws.Range("F34").Value = 5660
$ws.Range("F42").Value = 581
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 88
$ws.Range("F45").Value = 408
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 680
$ws.Range("F7").Value = 295
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 1219
$ws.Range("F7").Value = 680
$ws.Range("F8").Value = 74
$ws.Range("F10").Value = 295
$ws.Range("F15").Value = 618
$ws.Range("F24").Value = 1006
$ws.Range("C27").Value = "上海·坏孩纸物语第40届动漫节之曹沫篇"
$ws.Range("D27").Value = "沪南路2229号 复地活力城"
$ws.Range("E27").Value = "2024.05.02 13:00-05.02 17:40"
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 66.90000000000001
$ws.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=84724"
$ws.Range("I27").Value = "//i1.hdslb.com/bfs/openplatform/202404/Y4wTU9111713328435995.png"
$ws.Range("F37").Value = 799
$ws.Range("F38").Value = 144
$ws.Range("F39").Value = 1329
$ws.Range("F40").Value = 5660
$ws.Range("F48").Value = 581
$ws.Range("F52").Value = 408
